$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1094
$ws.Range("J19").Value = 1102.4546
$ws.Range("L19").Value = 1102.4546
$ws.Range("N19").Value = -1452.4546

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2772.5454
$ws.Range("I40").Value = 3166.5
$ws.Range("J40").Value = 2299.8
$ws.Range("K40").Value = 3166.5
$ws.Range("L40").Value = 2299.8
$ws.Range("M40").Value = -2991.5
$ws.Range("N40").Value = -2649.8

# ALC row 43
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 823880
$ws.Range("J43").Value = 1027875
$ws.Range("L43").Value = 1027875
$ws.Range("N43").Value = -1028013

# ALC row 48
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 1
$ws.Range("I48").Value = 1
$ws.Range("K48").Value = 3
$ws.Range("M48").Value = 289

# ALC row 56
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 1
$ws.Range("I56").Value = 1
$ws.Range("K56").Value = 3
$ws.Range("M56").Value = 531

# ALC row 62
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 125040210
$ws.Range("I62").Value = 250003740
$ws.Range("K62").Value = 250003740
$ws.Range("M62").Value = -250003116

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6733.3335
$ws.Range("I64").Value = 6500
$ws.Range("K64").Value = 6500
$ws.Range("M64").Value = -6252

# ALC row 65
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 125040210
$ws.Range("I65").Value = 250003740
$ws.Range("K65").Value = 1250018700
$ws.Range("M65").Value = -1250015580

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 6733.3335
$ws.Range("I67").Value = 6500
$ws.Range("K67").Value = 6500
$ws.Range("M67").Value = -5642

# ALC row 74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 62507470
$ws.Range("I74").Value = 100004750
$ws.Range("K74").Value = 100004750
$ws.Range("M74").Value = -100003814

# ALC row 77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 62507470
$ws.Range("I77").Value = 100004750
$ws.Range("K77").Value = 500023750
$ws.Range("M77").Value = -500019070

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 25704786
$ws.Range("I88").Value = 83336984
$ws.Range("J88").Value = 90477.22
$ws.Range("K88").Value = 83336984
$ws.Range("L88").Value = 90477.22
$ws.Range("M88").Value = -83336578
$ws.Range("N88").Value = -91289.22

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 25704786
$ws.Range("I91").Value = 83336984
$ws.Range("J91").Value = 90477.22
$ws.Range("K91").Value = 83336984
$ws.Range("L91").Value = 90477.22
$ws.Range("M91").Value = -83335580
$ws.Range("N91").Value = -93285.22

# ALC row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 6592.9062
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 6592.9062
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 19778.7186
$ws.Range("N112").Value = -21994.7186
$ws.Range("M112").ClearContents()

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1686.439
$ws.Range("I132").Value = 1461.2333
$ws.Range("J132").Value = 2300.6365
$ws.Range("K132").Value = 4383.699900000001
$ws.Range("L132").Value = 6901.9095
$ws.Range("M132").Value = -1853.699900000001
$ws.Range("N132").Value = -11961.9095

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2824.3572
$ws.Range("I137").Value = 3045.7
$ws.Range("K137").Value = 9137.099999999999
$ws.Range("M137").Value = -6587.099999999999

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3651.5103
$ws.Range("J138").Value = 5529.5864
$ws.Range("L138").Value = 16588.7592
$ws.Range("N138").Value = -26868.7592

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2137.0667
$ws.Range("I141").Value = 2150.077
$ws.Range("K141").Value = 6450.231000000001
$ws.Range("M141").Value = -1270.231000000001

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1999.25

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1999.25

# BSM row 35
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 59999.5
$ws.Range("J35").Value = 59999.5
$ws.Range("L35").Value = 59999.5
$ws.Range("N35").Value = -60619.5

# BSM row 40
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 53942
$ws.Range("J40").Value = 53942
$ws.Range("L40").Value = 53942
$ws.Range("N40").Value = -54472

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 139003870
$ws.Range("I86").Value = 41835820
$ws.Range("K86").Value = 41835820
$ws.Range("M86").Value = -41834697

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 139003870
$ws.Range("I89").Value = 41835820
$ws.Range("K89").Value = 209179100
$ws.Range("M89").Value = -209173484

# BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 6602.7144
$ws.Range("I94").Value = 604.5
$ws.Range("J94").Value = 9002
$ws.Range("K94").Value = 604.5
$ws.Range("L94").Value = 9002
$ws.Range("M94").Value = -153.5
$ws.Range("N94").Value = -9904

# BSM row 96
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 23967.143
$ws.Range("I96").Value = 11977.2
$ws.Range("J96").Value = 53942
$ws.Range("K96").Value = 11977.2
$ws.Range("L96").Value = 53942
$ws.Range("M96").Value = -9231.200000000001
$ws.Range("N96").Value = -59434

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5849.6
$ws.Range("I16").Value = 3335.111
$ws.Range("J16").Value = 7906.909
$ws.Range("K16").Value = 3335.111
$ws.Range("L16").Value = 7906.909
$ws.Range("M16").Value = -3048.111
$ws.Range("N16").Value = -8480.909

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8976.6
$ws.Range("I31").Value = 3802.5
$ws.Range("J31").Value = 10858.091
$ws.Range("K31").Value = 3802.5
$ws.Range("L31").Value = 10858.091
$ws.Range("M31").Value = -3507.5
$ws.Range("N31").Value = -11448.091

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 8976.6
$ws.Range("I34").Value = 3802.5
$ws.Range("J34").Value = 10858.091
$ws.Range("K34").Value = 3802.5
$ws.Range("L34").Value = 10858.091
$ws.Range("M34").Value = -3600.5
$ws.Range("N34").Value = -11262.091

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9621403
$ws.Range("I58").Value = 22730068
$ws.Range("J58").Value = 8381.733
$ws.Range("K58").Value = 22730068
$ws.Range("L58").Value = 8381.733
$ws.Range("M58").Value = -22729865
$ws.Range("N58").Value = -8787.733

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 3538.9285
$ws.Range("I105").Value = 1654.5
$ws.Range("J105").Value = 8250
$ws.Range("K105").Value = 1654.5
$ws.Range("L105").Value = 8250
$ws.Range("M105").Value = 92.5
$ws.Range("N105").Value = -11744

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 5849.6
$ws.Range("I113").Value = 3335.111
$ws.Range("J113").Value = 7906.909
$ws.Range("K113").Value = 3335.111
$ws.Range("L113").Value = 7906.909
$ws.Range("M113").Value = -1165.111
$ws.Range("N113").Value = -12246.909

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 7486.5093
$ws.Range("I134").Value = 6716.8965
$ws.Range("K134").Value = 20150.6895
$ws.Range("M134").Value = -17615.6895

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 9621403
$ws.Range("I136").Value = 22730068
$ws.Range("J136").Value = 8381.733
$ws.Range("K136").Value = 68190204
$ws.Range("L136").Value = 25145.199
$ws.Range("M136").Value = -68187654
$ws.Range("N136").Value = -30245.199

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2002367.6
$ws.Range("I5").Value = 3334558.5
$ws.Range("J5").Value = 4081.25
$ws.Range("K5").Value = 10003675.5
$ws.Range("L5").Value = 12243.75
$ws.Range("M5").Value = -10003563.5
$ws.Range("N5").Value = -12467.75

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1802.0488
$ws.Range("J131").Value = 2033.4
$ws.Range("L131").Value = 6100.200000000001
$ws.Range("N131").Value = -16180.2

# CUL row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 13635.667
$ws.Range("I132").Value = 6000
$ws.Range("J132").Value = 19744.2
$ws.Range("K132").Value = 54000
$ws.Range("L132").Value = 177697.8
$ws.Range("M132").Value = -51470
$ws.Range("N132").Value = -182757.8

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 2002367.6
$ws.Range("I135").Value = 3334558.5
$ws.Range("J135").Value = 4081.25
$ws.Range("K135").Value = 30011026.5
$ws.Range("L135").Value = 36731.25
$ws.Range("M135").Value = -30008491.5
$ws.Range("N135").Value = -41801.25

# CUL row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 334164.5
$ws.Range("I137").Value = 250997
$ws.Range("K137").Value = 752991
$ws.Range("M137").Value = -747891

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2095.182
$ws.Range("I97").Value = 1594.7826
$ws.Range("J97").Value = 3246.1
$ws.Range("K97").Value = 1594.7826
$ws.Range("L97").Value = 3246.1
$ws.Range("M97").Value = -1098.7826
$ws.Range("N97").Value = -4238.1

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 343331
$ws.Range("I122").Value = 507497
$ws.Range("J122").Value = 14999
$ws.Range("K122").Value = 1522491
$ws.Range("L122").Value = 44997
$ws.Range("M122").Value = -1520041
$ws.Range("N122").Value = -49897

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1569.5
$ws.Range("I16").Value = 1458.6666
$ws.Range("K16").Value = 1458.6666
$ws.Range("M16").Value = -1288.6666

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 10102804
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 10102804
$ws.Range("N46").Value = -10103180
$ws.Range("M46").ClearContents()

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 125000750
$ws.Range("I55").Value = 333333540
$ws.Range("J55").Value = 1080
$ws.Range("K55").Value = 333333540
$ws.Range("L55").Value = 1080
$ws.Range("M55").Value = -333333367
$ws.Range("N55").Value = -1426

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4218.8276
$ws.Range("I122").Value = 3111.476
$ws.Range("K122").Value = 9334.428
$ws.Range("M122").Value = -6884.428

# LTW row 127
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 59519
$ws.Range("J127").Value = 59519
$ws.Range("L127").Value = 59519
$ws.Range("N127").Value = -69439

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 37038176
$ws.Range("I107").Value = 1149.4
$ws.Range("K107").Value = 3448.2
$ws.Range("M107").Value = -1528.2

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 290871.66
$ws.Range("I122").Value = 1336401
$ws.Range("K122").Value = 4009203
$ws.Range("M122").Value = -4006753

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4700.2856
$ws.Range("I126").Value = 3380.4
$ws.Range("K126").Value = 10141.2
$ws.Range("M126").Value = -7671.200000000001

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 58885836
$ws.Range("I136").Value = 125000744
$ws.Range("K136").Value = 375002232
$ws.Range("M136").Value = -374999682
